# Remove the "maintenance_log_image" survey question row (the image-capture
# field) from the "survey" sheet. This is the 10th row of that sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")
$ws.Rows.Item(10).Delete()
